$d = $word.ActiveDocument

# 1) Re-number the "smartrewards" bookmark id (1 -> 0) by deleting and re-adding it
$bm = $d.Bookmarks(1)
$bmRange = $bm.Range
$bmName = $bm.Name
$bm.Delete()
$d.Bookmarks.Add($bmName, $bmRange) | Out-Null

# 2) Replace the SmartRewards intro paragraph: simplify its paragraph
#    properties and split the long run that mentions the Decentralized
#    Distribution change into three runs (one with new wording/formatting).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "SmartRewards*are a price stabilization*") {
        $target = $p
        break
    }
}

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>SmartRewards</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">&#x00a0;are a price stabilization mechanism and a way to encourage long term holding. Long term holders are key to the project’s success since the SmartHive treasury needs SmartCash to appreciate in value in order to fund meaningful 3rd party proposals and help grow SmartCash into a successful global crypto-currency. Beginning at block 574,100, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">the current SmartRewards model changed to a Decentralized Distribution</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">. SmartRewards will then be distributed after every&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">47,500 Blocks</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">&#x00a0;to all eligible addresses. The payouts will begin 200 Blocks after the cycle ends and 1000 addresses will be paid every other block. SmartNodes will keep track on which addresses are eligible and which are not. You will earn SmartRewards on&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">&#x00a0;address for which you hold the keys (web or desktop, SmartNode included!) which holds &gt;=1000 SMART for one month and does not make&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">outgoing transactions during that time. Please note, most exchanges do&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:i/><w:iCs/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">&#x00a0;pay SmartRewards to their users, holding &gt;= 1000 SMART on an exchange does not guarantee a reward</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">.&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">The SmartRewards will come out of the 15% block reward allocation.</w:t></w:r></w:p>'
$target.Range.InsertXML($newParaXml)

Write-Output "done"
